# "add unittests for whitespace sanitizer"
#
# Adds a new worksheet, LAST_COL_WHITESPACE, right after Sheet1. It holds
# the same HEADER_A / whitespace-before / whitespace-after data as
# Sheet1's columns A and B, but (deliberately, for the "last column is
# blank/whitespace" test case) without a header in B1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LAST_COL_WHITESPACE"

# Copy column A (all rows) and column B (rows 2-15, skipping the header
# row so B1 stays empty) from Sheet1 into the new sheet.
for ($r = 1; $r -le 15; $r++) {
  $ws2.Cells.Item($r, 1).Value2 = $ws1.Cells.Item($r, 1).Value()
  if ($r -ge 2) {
    $ws2.Cells.Item($r, 2).Value2 = $ws1.Cells.Item($r, 2).Value()
  }
}

# Match column B's width on the new sheet.
$ws2.Columns.Item(2).ColumnWidth = 20.5

# Sheet1 is no longer the active tab; update its selection.
$ws1.Activate()
$ws1.Range("A1:B15").Select()

# Make the new sheet the active tab, with cell B1 selected.
$ws2.Activate()
$ws2.Range("B1").Select()
